$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.564
$ws.Range("A6").Value = -22.255
$ws.Range("A7").Value = -19.575
$ws.Range("C7").Value = -12.489
$ws.Range("C12").Value = -11.06
$ws.Range("E13").Value = 16.641
$ws.Range("E14").Value = 17.13
$ws.Range("C15").Value = -13.447
$ws.Range("A16").Value = -21.96
$ws.Range("E16").Value = 16.683
$ws.Range("E19").Value = 16.45
$ws.Range("A20").Value = -20.026
$ws.Range("C20").Value = -12.245
$ws.Range("C21").Value = -12.519
$ws.Range("C22").Value = -12.624
$ws.Range("E22").Value = 16.82
$ws.Range("C23").Value = -12.223
$ws.Range("A28").Value = -21.942
$ws.Range("A29").Value = -21.524
$ws.Range("C29").Value = -12.163
$ws.Range("A32").Value = -21.864
$ws.Range("C34").Value = -11.955
$ws.Range("E36").Value = 16.709
$ws.Range("A40").Value = -20.097
$ws.Range("C42").Value = -12.492
$ws.Range("C43").Value = -13.159
$ws.Range("C44").Value = -13.216
$ws.Range("C45").Value = -13.054
$ws.Range("A46").Value = -21.995
$ws.Range("C46").Value = -13.646
$ws.Range("E46").Value = 16.835
$ws.Range("C50").Value = -14.005
$ws.Range("E50").Value = 16.434
$ws.Range("A51").Value = -21.591
$ws.Range("C51").Value = -11.06
$ws.Range("A52").Value = -21.872
$ws.Range("A57").Value = -22.247
$ws.Range("A59").Value = -22.405
$ws.Range("A62").Value = -22.158
$ws.Range("A66").Value = -21.591
$ws.Range("C66").Value = -11.405
$ws.Range("C67").Value = -11.354
$ws.Range("A73").Value = -20.597
$ws.Range("A74").Value = -21.244
$ws.Range("C79").Value = -11.816
$ws.Range("C84").Value = -14.098
$ws.Range("A92").Value = -21.609
$ws.Range("C92").Value = -11.383
$ws.Range("E95").Value = 17.385
$ws.Range("C97").Value = -11.87
$ws.Range("E97").Value = 17.197
$ws.Range("A100").Value = -22.063
